# Update intro section of project
# 1. Change author-line affiliation superscripts "3" -> "1" (Steven Cognac and
#    Alex Clippinger both reference affiliation 1, not 3).
# 2. Remove the now-unused third affiliation line ("3 UCSB") from the
#    affiliations paragraph.
# 3. Update the "report was generated on" timestamp.
# 4. Append the actual Git commit details block under the Colophon section.

$d = $word.ActiveDocument

# --- 1. Superscript "3" -> "1" on the two author-name references ------------
# The author-byline paragraph contains exactly two "3" characters (both
# superscripted footnote markers, one after each author name). Scoping the
# Find/Replace to this paragraph's Range and asking for "replace all" (the
# trailing 2 == wdReplaceAll) swaps both of them in place, preserving the
# surrounding run formatting (vertAlign superscript).
$authorRange = $d.Paragraphs(2).Range
$authorRange.Find.ClearFormatting()
$authorRange.Find.Execute("3", $true, $false, $false, $false, $false, $true, 1, $false, "1", 2)

# --- 2. Drop the third affiliation line ("3 UCSB") ---------------------------
# The affiliations paragraph currently reads:
#   "1 University of California Santa Barbara" <br>
#   "2 University of California Santa Barbara" <br>
#   "3 UCSB"
# Remove the trailing <br> + "3" + " " + "UCSB" portion, leaving only the
# first two affiliation lines.
$affilPara = $d.Paragraphs(5)
$affilText = $affilPara.Range.Text
$marker = "Santa Barbara"
$lastIdx = $affilText.LastIndexOf($marker)
$cutStart = $affilPara.Range.Start + $lastIdx + $marker.Length
$cutEnd = $affilPara.Range.End - 1
$cutRange = $d.Range($cutStart, $cutEnd)
$cutRange.Delete()

# --- 3. Update the generated-on timestamp ------------------------------------
$d.Content.Find.Execute("This report was generated on 2021-08-24 11:42:26 using", $true, $false, $false, $false, $false, $true, 1, $false, "This report was generated on 2021-08-24 11:53:13 using", 2)

# --- 4. Append the Git commit details block ----------------------------------
$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$lastPara.Range.InsertParagraphAfter()
$gitPara = $d.Paragraphs($d.Paragraphs.Count)
$gitPara.Style = "Source Code"

$line1 = "#> Local:    master C:/Users/stvec/Documents/mypaper"
$line2 = "#> Remote:   master @ origin (https://github.com/cognack/Stevens_paper.git)"
$line3 = "#> Head:     [a97e68a] 2021-08-24: Add GitHub links to DESCRIPTION"

$insertPos = $gitPara.Range.End - 1
$cur = $d.Range($insertPos, $insertPos)
$cur.InsertAfter($line1)
$t1 = $d.Range($insertPos, $insertPos + $line1.Length)
$t1.Style = "Verbatim Char"

$insertPos = $gitPara.Range.End - 1
$cur = $d.Range($insertPos, $insertPos)
$cur.InsertBreak(6)

$insertPos = $gitPara.Range.End - 1
$cur = $d.Range($insertPos, $insertPos)
$cur.InsertAfter($line2)
$t2 = $d.Range($insertPos, $insertPos + $line2.Length)
$t2.Style = "Verbatim Char"

$insertPos = $gitPara.Range.End - 1
$cur = $d.Range($insertPos, $insertPos)
$cur.InsertBreak(6)

$insertPos = $gitPara.Range.End - 1
$cur = $d.Range($insertPos, $insertPos)
$cur.InsertAfter($line3)
$t3 = $d.Range($insertPos, $insertPos + $line3.Length)
$t3.Style = "Verbatim Char"

Write-Output "Edit complete"
